$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reallocate Director for Tier2 (row 3) from David to Conrad
$ws.Range("C3").Value = "Conrad"

# Update the active selection (mirrors the saved selection state in the file)
$ws.Range("E8").Select()
